# Update final evaluation results across all three sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.6226591760299626
$wsSummary.Range("C2").Value = 0.5771495877502945
$wsSummary.Range("D2").Value = 0.9176029962546817
$wsSummary.Range("E2").Value = 0.7086044830079538
$wsSummary.Range("F2").Value = 0.8207705192629816
$wsSummary.Range("G2").Value = 0.8972462849496443
$wsSummary.Range("H2").Value = 0.7754176661195977
$wsSummary.Range("I2").Value = 490
$wsSummary.Range("J2").Value = 359
$wsSummary.Range("K2").Value = 175
$wsSummary.Range("L2").Value = 44

# --- Sheet: Classification Report ---
$wsReport = $wb.Worksheets.Item("Classification Report")
$wsReport.Range("B2").Value = 0.7990867579908676
$wsReport.Range("C2").Value = 0.3277153558052435
$wsReport.Range("D2").Value = 0.4648074369189907

$wsReport.Range("B3").Value = 0.5771495877502945
$wsReport.Range("C3").Value = 0.9176029962546817
$wsReport.Range("D3").Value = 0.7086044830079538

$wsReport.Range("B4").Value = 0.6226591760299626
$wsReport.Range("C4").Value = 0.6226591760299626
$wsReport.Range("D4").Value = 0.6226591760299626
$wsReport.Range("E4").Value = 0.6226591760299626

$wsReport.Range("B5").Value = 0.688118172870581
$wsReport.Range("C5").Value = 0.6226591760299626
$wsReport.Range("D5").Value = 0.5867059599634722

$wsReport.Range("B6").Value = 0.6881181728705811
$wsReport.Range("C6").Value = 0.6226591760299626
$wsReport.Range("D6").Value = 0.5867059599634722

# --- Sheet: Confusion Matrix ---
$wsConfusion = $wb.Worksheets.Item("Confusion Matrix")
$wsConfusion.Range("B2").Value = 175
$wsConfusion.Range("C2").Value = 359
$wsConfusion.Range("B3").Value = 44
$wsConfusion.Range("C3").Value = 490
